$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales")

# --- Row 77: update Reference (D77) and ensure existing columns remain intact ---
$ws.Range('D77').Value = 'Invoice#348, QB#5052'

# --- Row 84: remove the stray empty Reference/Actions placeholder in column H ---
$ws.Range('H84').ClearContents()

# --- Append new transaction rows 120-128 ---
# Row 120
$ws.Range('A120').Value = 'علي عمر دين'
$ws.Range('B120').NumberFormat = '@'
$ws.Range('B120').Value = '2025-07-29'
$ws.Range('C120').Value = '#1: أكريليك استاند إي فور مع ثنية | Qty: 1 | Price: 35 | Total: 35 | VAT: 0'
$ws.Range('E120').Value = 35
$ws.Range('F120').Value = 0
$ws.Range('G120').Value = 35
$ws.Range('I120').Value = $false

# Row 121
$ws.Range('A121').Value = 'علا للمصاعد'
$ws.Range('B121').NumberFormat = '@'
$ws.Range('B121').Value = '2025-07-29'
$ws.Range('C121').Value = '#1: سند صيانة | Qty: 50 | Price: 13 | Total: 650 | VAT: 97.5'
$ws.Range('E121').Value = 650
$ws.Range('F121').Value = 97.5
$ws.Range('G121').Value = 747.5
$ws.Range('I121').Value = $false

# Row 122
$ws.Range('A122').Value = 'رائد - كابلات بحرة'
$ws.Range('B122').NumberFormat = '@'
$ws.Range('B122').Value = '2025-07-30'
$ws.Range('C122').Value = '#1: ختم بيضاوي 3045 | Qty: 1 | Price: 140 | Total: 140 | VAT: 0'
$ws.Range('D122').Value = 'QB#5083'
$ws.Range('E122').Value = 140
$ws.Range('F122').Value = 0
$ws.Range('G122').Value = 140
$ws.Range('I122').Value = $false

# Row 123
$ws.Range('A123').Value = 'شاي فال'
$ws.Range('B123').NumberFormat = '@'
$ws.Range('B123').Value = '2025-07-30'
$ws.Range('C123').Value = '#1: منيو بوخار | Qty: 10 | Price: 10 | Total: 100 | VAT: 15; #2: منيو شاي فال | Qty: 10 | Price: 2.75 | Total: 27.5 | VAT: 4.13; #3: كرت افتتاح فرع 9.5*13.5 سم | Qty: 500 | Price: 0.6 | Total: 300 | VAT: 45; #4: أكريليك مقاس 35*60 سم | Qty: 2 | Price: 250 | Total: 500 | VAT: 75'
$ws.Range('D123').Value = 'QB#50845085'
$ws.Range('E123').Value = 927.5
$ws.Range('F123').Value = 139.13
$ws.Range('G123').Value = 1066.63
$ws.Range('I123').Value = $false

# Row 124
$ws.Range('A124').Value = 'خالد أبو سعيد'
$ws.Range('B124').NumberFormat = '@'
$ws.Range('B124').Value = '2025-07-30'
$ws.Range('C124').Value = '#1: استيكرات قهوة متنوع | Qty: 1 | Price: 0 | Total: 0 | VAT: 0'
$ws.Range('E124').Value = 0
$ws.Range('F124').Value = 0
$ws.Range('G124').Value = 0
$ws.Range('I124').Value = $false

# Row 125
$ws.Range('A125').Value = 'محمد المالكي UK صاحب صهيب'
$ws.Range('B125').NumberFormat = '@'
$ws.Range('B125').Value = '2025-07-29'
$ws.Range('C125').Value = '#1: طباعة كيس ورقي للهدايا | Qty: 3 | Price: 25 | Total: 75 | VAT: 0'
$ws.Range('E125').Value = 75
$ws.Range('F125').Value = 0
$ws.Range('G125').Value = 75
$ws.Range('I125').Value = $false

# Row 126
$ws.Range('A126').Value = 'عبدالله السندي - صاحب عمر'
$ws.Range('B126').NumberFormat = '@'
$ws.Range('B126').Value = '2025-07-29'
$ws.Range('C126').Value = '#1: طباعة ورق إي ثري ديجيتال - بلاش | Qty: 3 | Price: 0 | Total: 0 | VAT: 0; #2: بنر مقاس 50*70 سم | Qty: 1 | Price: 40 | Total: 40 | VAT: 0; #3: بروشور 150 جرام A5 طباعة وجهين | Qty: 100 | Price: 1.2 | Total: 120 | VAT: 0; #4: فلين مقاس 100*70 سم | Qty: 2 | Price: 75 | Total: 150 | VAT: 0; #5: فلين مقاس 50*50 سم | Qty: 5 | Price: 40 | Total: 200 | VAT: 0; #6: خصم 10 ريال | Qty: 1 | Price: 0 | Total: 0 | VAT: 0'
$ws.Range('E126').Value = 510
$ws.Range('F126').Value = 0
$ws.Range('G126').Value = 510
$ws.Range('I126').Value = $false

# Row 127
$ws.Range('A127').Value = 'البروج الذهبية'
$ws.Range('B127').NumberFormat = '@'
$ws.Range('B127').Value = '2025-07-19'
$ws.Range('C127').Value = '#1: تصميم الهوية | Qty: 1 | Price: 304.35 | Total: 304.35 | VAT: 45.65; #2: طباعة فولدر مقاس A4 - طباعة جهة واحدة - بجيب داخلي واحد على اليسار | Qty: 100 | Price: 5 | Total: 500 | VAT: 75; #3: طباعة ورق الخطابات كونكورر فاخر | Qty: 500 | Price: .8 | Total: 400 | VAT: 60; #4: طباعة ورق الخطابات وودفري | Qty: 500 | Price: .6 | Total: 300 | VAT: 45; #5: طباعة ظرف A4 وجه واحد | Qty: 200 | Price: 2.6 | Total: 520 | VAT: 78; #6: طباعة ظرف A5 وجه واحد | Qty: 200 | Price: 2.2 | Total: 440.00000000000006 | VAT: 66; #7: طباعة ظرف DL وجه واحد | Qty: 200 | Price: 1.8 | Total: 360 | VAT: 54; #8: طباعة سندات قبض - لون واحد - 1+2 | Qty: 10 | Price: 18 | Total: 180 | VAT: 27; #9: طباعة سندات صرف - لون واحد - 1+2 | Qty: 10 | Price: 18 | Total: 180 | VAT: 27; #10: ختم دائري R538 | Qty: 1 | Price: 120 | Total: 120 | VAT: 18'
$ws.Range('D127').Value = 'Quotation#194, Invoice#347'
$ws.Range('E127').Value = 3304.35
$ws.Range('F127').Value = 495.65
$ws.Range('G127').Value = 3800
$ws.Range('I127').Value = $false

# Row 128
$ws.Range('A128').Value = 'شركة كيري'
$ws.Range('B128').NumberFormat = '@'
$ws.Range('B128').Value = '2025-07-30'
$ws.Range('C128').Value = '#1: ختم دائري R538 | Qty: 2 | Price: 120 | Total: 240 | VAT: 36'
$ws.Range('E128').Value = 240
$ws.Range('F128').Value = 36
$ws.Range('G128').Value = 276
$ws.Range('I128').Value = $false
